# Apply the "Dataset Loading" changes to the diagnostics sheet:
#  1. Fix the G102 timestamp (rounding correction).
#  2. Append 10 new log rows (103-112) for AIDS / PTC_FR / MUTAG / Letter-high
#     SVC_Simple_Prototype_GED_poly functionality-test runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the existing G102 timestamp value -------------------------
$ws.Range("G102").Value2 = 45905.77089622685

# --- 2. Append the new rows -------------------------------------------------
# Template number format copied from an existing timestamp cell so the new
# F/G cells line up with the same style (YYYY-MM-DD HH:MM:SS) used elsewhere.
$tsFormat = $ws.Range("G2").NumberFormat

$rows = @(
    @{ Row=103; A="Fucntionality_test_AIDS_with_SVC_Simple_Prototype_GED_poly";        B="AIDS";        C="SVC_Simple_Prototype_GED_poly"; D="SVC_Simple_Prototype_GED_poly_trained_on_AIDS.joblib";        F=45910.64188774306; G=45910.64188697917 },
    @{ Row=104; A="Fucntionality_test_AIDS_with_SVC_Simple_Prototype_GED_poly";        B="AIDS";        C="SVC_Simple_Prototype_GED_poly"; D="SVC_Simple_Prototype_GED_poly_trained_on_AIDS.joblib";        F=45910.64188774306; G=45910.64188697917 },
    @{ Row=105; A="Fucntionality_test_PTC_FR_with_SVC_Simple_Prototype_GED_poly";      B="PTC_FR";      C="SVC_Simple_Prototype_GED_poly"; D="SVC_Simple_Prototype_GED_poly_trained_on_PTC_FR.joblib";      F=45910.6429222338;  G=45910.64292206019 },
    @{ Row=106; A="Fucntionality_test_PTC_FR_with_SVC_Simple_Prototype_GED_poly";      B="PTC_FR";      C="SVC_Simple_Prototype_GED_poly"; D="SVC_Simple_Prototype_GED_poly_trained_on_PTC_FR.joblib";      F=45910.6429222338;  G=45910.64292206019 },
    @{ Row=107; A="Fucntionality_test_MUTAG_with_SVC_Simple_Prototype_GED_poly";       B="MUTAG";       C="SVC_Simple_Prototype_GED_poly"; D="SVC_Simple_Prototype_GED_poly_trained_on_MUTAG.joblib";       F=45910.64468162037; G=45910.64468146991 },
    @{ Row=108; A="Fucntionality_test_MUTAG_with_SVC_Simple_Prototype_GED_poly";       B="MUTAG";       C="SVC_Simple_Prototype_GED_poly"; D="SVC_Simple_Prototype_GED_poly_trained_on_MUTAG.joblib";       F=45910.64468162037; G=45910.64468146991 },
    @{ Row=109; A="Fucntionality_test_MUTAG_with_SVC_Simple_Prototype_GED_poly";       B="MUTAG";       C="SVC_Simple_Prototype_GED_poly"; D="SVC_Simple_Prototype_GED_poly_trained_on_MUTAG.joblib";       F=45910.64681074074; G=45910.64681050926 },
    @{ Row=110; A="Fucntionality_test_MUTAG_with_SVC_Simple_Prototype_GED_poly";       B="MUTAG";       C="SVC_Simple_Prototype_GED_poly"; D="SVC_Simple_Prototype_GED_poly_trained_on_MUTAG.joblib";       F=45910.64681074074; G=45910.64681050926 },
    @{ Row=111; A="Fucntionality_test_Letter-high_with_SVC_Simple_Prototype_GED_poly"; B="Letter-high"; C="SVC_Simple_Prototype_GED_poly"; D="SVC_Simple_Prototype_GED_poly_trained_on_Letter-high.joblib"; F=45910.65059689815; G=45910.65059133102 },
    @{ Row=112; A="Fucntionality_test_Letter-high_with_SVC_Simple_Prototype_GED_poly"; B="Letter-high"; C="SVC_Simple_Prototype_GED_poly"; D="SVC_Simple_Prototype_GED_poly_trained_on_Letter-high.joblib"; F=45910.65059689566; G=45910.65059132903 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D

    # E (dataset_load_duration) is present but blank in the source log, same
    # as the rest of the sheet - write it (as a lone quote/text prefix, the
    # classic "force empty text" trick) then strip the style back to Normal
    # so the cell round-trips as an empty text cell instead of a disappearing
    # blank one.
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = "'"
    $eCell.Style = "Normal"

    $fCell = $ws.Cells.Item($row, 6)
    $fCell.NumberFormat = $tsFormat
    $fCell.Value2 = $r.F

    $gCell = $ws.Cells.Item($row, 7)
    $gCell.NumberFormat = $tsFormat
    $gCell.Value2 = $r.G

    # H (Error) is likewise an empty text cell.
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value = "'"
    $hCell.Style = "Normal"
}
